$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Student ID column (A) values for rows 2-13.
# Leading apostrophe keeps these as text (matching the original inlineStr
# cells), since the Student ID column holds text-typed numeric codes.
$ws.Range("A2").Value = "'190874"
$ws.Range("A3").Value = "'201253"
$ws.Range("A4").Value = "'201252"
$ws.Range("A5").Value = "'201023"
$ws.Range("A6").Value = "'201670"
$ws.Range("A7").Value = "'190796"
$ws.Range("A8").Value = "'201838"
$ws.Range("A9").Value = "'191258"
$ws.Range("A10").Value = "'200468"
$ws.Range("A11").Value = "'201065"
$ws.Range("A12").Value = "'191502"
$ws.Range("A13").Value = "'200933"

# Delete rows 14-25 (old extra data no longer present)
$ws.Rows("14:25").Delete()
